# Generate Report for Handback
# Both locales (zh-cn, de-de) are now in sync with the en-US source, so the
# status flips from "Ready for handoff" to "Handed back: in sync with en-US",
# the handback timestamps advance, and the stale-handback-version warning
# clears since the handback file is current.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: per-locale status columns ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus
$ovw.Columns.Item(5).AutoFit()
$ovw.Columns.Item(6).AutoFit()

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("K2").Value = "2016-08-15 22:46:31"
$zh.Range("P2").Value = ""
$zh.Columns.Item(3).AutoFit()
$zh.Columns.Item(16).AutoFit()

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("K2").Value = "2016-08-15 22:46:38"
$de.Range("P2").Value = ""
$de.Columns.Item(3).AutoFit()
$de.Columns.Item(16).AutoFit()
